$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("addBrandNewVisitor")

$ws.Range("A2").Value = "SeleniumTest+v20200128122914@gmail.com"
$ws.Range("B2").Value = "Test+v20200128122914"
